$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 657.3333
$ws.Range("I38").Value = 121.6
$ws.Range("J38").Value = 1040
$ws.Range("K38").Value = 364.8
$ws.Range("L38").Value = 3120
$ws.Range("M38").Value = 7.200000000000045
$ws.Range("N38").Value = -3864
$ws.Range("H39").Value = 153.4
$ws.Range("I39").Value = 98
$ws.Range("J39").Value = 375
$ws.Range("K39").Value = 294
$ws.Range("L39").Value = 1125
$ws.Range("M39").Value = 2
$ws.Range("N39").Value = -1717
$ws.Range("H64").Value = 3957.9167
$ws.Range("I64").Value = 3473.75
$ws.Range("J64").Value = 4200
$ws.Range("K64").Value = 3473.75
$ws.Range("L64").Value = 4200
$ws.Range("M64").Value = -3225.75
$ws.Range("N64").Value = -4696
$ws.Range("H67").Value = 3957.9167
$ws.Range("I67").Value = 3473.75
$ws.Range("J67").Value = 4200
$ws.Range("K67").Value = 3473.75
$ws.Range("L67").Value = 4200
$ws.Range("M67").Value = -2615.75
$ws.Range("N67").Value = -5916
$ws.Range("H74").Value = 4138.143
$ws.Range("I74").Value = 4292.385
$ws.Range("K74").Value = 4292.385
$ws.Range("M74").Value = -3356.385
$ws.Range("H76").Value = 3667.3333
$ws.Range("I76").Value = 3200
$ws.Range("J76").Value = 3901
$ws.Range("K76").Value = 3200
$ws.Range("L76").Value = 3901
$ws.Range("M76").Value = -2885
$ws.Range("N76").Value = -4531
$ws.Range("H77").Value = 4138.143
$ws.Range("I77").Value = 4292.385
$ws.Range("K77").Value = 21461.925
$ws.Range("M77").Value = -16781.925
$ws.Range("H79").Value = 3667.3333
$ws.Range("I79").Value = 3200
$ws.Range("J79").Value = 3901
$ws.Range("K79").Value = 3200
$ws.Range("L79").Value = 3901
$ws.Range("M79").Value = -2108
$ws.Range("N79").Value = -6085
$ws.Range("H113").Value = 3251.037
$ws.Range("I113").Value = 2888.3157
$ws.Range("K113").Value = 2888.3157
$ws.Range("M113").Value = 365.6842999999999
$ws.Range("H121").Value = 1256.5714
$ws.Range("J121").Value = 1315.1578
$ws.Range("L121").Value = 3945.4734
$ws.Range("N121").Value = -7439.4734

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H38").Value = 10000
$ws.Range("J38").Value = 10000
$ws.Range("L38").Value = 10000
$ws.Range("N38").Value = -10934
$ws.Range("H54").Value = 12210
$ws.Range("J54").Value = 12210
$ws.Range("L54").Value = 12210
$ws.Range("N54").Value = -13748
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1119.0834
$ws.Range("I64").Value = 616
$ws.Range("J64").Value = 1622.1666
$ws.Range("K64").Value = 616
$ws.Range("L64").Value = 1622.1666
$ws.Range("M64").Value = -391
$ws.Range("N64").Value = -2072.1666
$ws.Range("H67").Value = 1119.0834
$ws.Range("I67").Value = 616
$ws.Range("J67").Value = 1622.1666
$ws.Range("K67").Value = 616
$ws.Range("L67").Value = 1622.1666
$ws.Range("M67").Value = 164
$ws.Range("N67").Value = -3182.1666
$ws.Range("H134").Value = 2020.5428
$ws.Range("I134").Value = 1203.0769
$ws.Range("J134").Value = 4382.1113
$ws.Range("K134").Value = 3609.2307
$ws.Range("L134").Value = 13146.3339
$ws.Range("M134").Value = -1074.2307
$ws.Range("N134").Value = -18216.3339
$ws.Range("H137").Value = 49780
$ws.Range("J137").Value = 49780
$ws.Range("L137").Value = 49780
$ws.Range("N137").Value = -59980
$ws.Range("H138").Value = 24797.143
$ws.Range("J138").Value = 24797.143
$ws.Range("L138").Value = 24797.143
$ws.Range("N138").Value = -35077.143
$ws.Range("H140").Value = 54020
$ws.Range("J140").Value = 54020
$ws.Range("L140").Value = 54020
$ws.Range("N140").Value = -64380

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 5014000
$ws.Range("J4").Value = 5014000
$ws.Range("L4").Value = 5014000
$ws.Range("N4").Value = -5014224
$ws.Range("H38").Value = 10042
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 10042
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 10042
$ws.Range("N38").Value = -10796
$ws.Range("H46").Value = 10042
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 10042
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 10042
$ws.Range("N46").Value = -10464
$ws.Range("H58").Value = 17242654
$ws.Range("I58").Value = 19231922
$ws.Range("J58").Value = 2330.3333
$ws.Range("K58").Value = 19231922
$ws.Range("L58").Value = 2330.3333
$ws.Range("M58").Value = -19231719
$ws.Range("N58").Value = -2736.3333
$ws.Range("H105").Value = 1789.9445
$ws.Range("I105").Value = 1729.1428
$ws.Range("J105").Value = 2002.75
$ws.Range("K105").Value = 1729.1428
$ws.Range("L105").Value = 2002.75
$ws.Range("M105").Value = 17.85719999999992
$ws.Range("N105").Value = -5496.75
$ws.Range("H132").Value = 39922.098
$ws.Range("I132").Value = 24678.35
$ws.Range("J132").Value = 112753.336
$ws.Range("K132").Value = 74035.04999999999
$ws.Range("L132").Value = 338260.008
$ws.Range("M132").Value = -71505.04999999999
$ws.Range("N132").Value = -343320.008
$ws.Range("H134").Value = 31504.916
$ws.Range("I134").Value = 1517.2693
$ws.Range("J134").Value = 109472.8
$ws.Range("K134").Value = 4551.8079
$ws.Range("L134").Value = 328418.4
$ws.Range("M134").Value = -2016.8079
$ws.Range("N134").Value = -333488.4
$ws.Range("H135").Value = 40000
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 40000
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 40000
$ws.Range("N135").Value = -50140
$ws.Range("H136").Value = 17242654
$ws.Range("I136").Value = 19231922
$ws.Range("J136").Value = 2330.3333
$ws.Range("K136").Value = 57695766
$ws.Range("L136").Value = 6990.999899999999
$ws.Range("M136").Value = -57693216
$ws.Range("N136").Value = -12090.9999
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("H138").Value = 47990
$ws.Range("J138").Value = 47990
$ws.Range("L138").Value = 47990
$ws.Range("N138").Value = -58270
$ws.Range("H140").Value = 51000
$ws.Range("J140").Value = 51000
$ws.Range("L140").Value = 51000
$ws.Range("N140").Value = -61360
$ws.Range("M38").ClearContents()
$ws.Range("M46").ClearContents()
$ws.Range("M135").ClearContents()
$ws.Range("N137").ClearContents()

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1893.4375
$ws.Range("I132").Value = 1648.75
$ws.Range("J132").Value = 2138.125
$ws.Range("K132").Value = 14838.75
$ws.Range("L132").Value = 19243.125
$ws.Range("M132").Value = -12308.75
$ws.Range("N132").Value = -24303.125
$ws.Range("H140").Value = 2483.581
$ws.Range("J140").Value = 2248.5576
$ws.Range("L140").Value = 6745.6728
$ws.Range("N140").Value = -17105.6728

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 2850
$ws.Range("J33").Value = 2850
$ws.Range("L33").Value = 2850
$ws.Range("N33").Value = -3354
$ws.Range("H62").Value = 25000
$ws.Range("I62").Value = 25000
$ws.Range("J62").Value = 25000
$ws.Range("K62").Value = 25000
$ws.Range("L62").Value = 25000
$ws.Range("M62").Value = -24314
$ws.Range("N62").Value = -26372
$ws.Range("H65").Value = 25000
$ws.Range("I65").Value = 25000
$ws.Range("J65").Value = 25000
$ws.Range("K65").Value = 75000
$ws.Range("L65").Value = 75000
$ws.Range("M65").Value = -71568
$ws.Range("N65").Value = -81864

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 5000
$ws.Range("J5").Value = 5000
$ws.Range("L5").Value = 5000
$ws.Range("N5").Value = -5224
